$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.420981
$ws.Range("H2").Value = 1.262943
$ws.Range("I2").Value = 0.005858863598414047
$ws.Range("J2").Value = 0.005858863598414048
$ws.Range("M2").Value = 1.937269333333333
$ws.Range("N2").Value = 5.811808
$ws.Range("O2").Value = 0.5832046952539398
$ws.Range("P2").Value = 0.5832046952539398
$ws.Range("Q2").Value = 0.815553581216
$ws.Range("R2").Value = 7.339982230944
$ws.Range("S2").Value = 0.003416916759447466
$ws.Range("T2").Value = 0.003416916759447466
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.420981
$ws.Range("H3").Value = 1.262943
$ws.Range("I3").Value = 0.005858863598414047
$ws.Range("J3").Value = 0.005858863598414048
$ws.Range("M3").Value = 0.9034129999999999
$ws.Range("N3").Value = 2.710239
$ws.Range("O3").Value = 0.2719677095424251
$ws.Range("P3").Value = 0.2719677095424251
$ws.Range("Q3").Value = 0.380319708153
$ws.Range("R3").Value = 3.422877373376999
$ws.Range("S3").Value = 0.001593421713382159
$ws.Range("T3").Value = 0.001593421713382159
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.420981
$ws.Range("H4").Value = 1.262943
$ws.Range("I4").Value = 0.005858863598414047
$ws.Range("J4").Value = 0.005858863598414048
$ws.Range("M4").Value = 0.4810833333333333
$ws.Range("N4").Value = 1.44325
$ws.Range("O4").Value = 0.1448275952036352
$ws.Range("P4").Value = 0.1448275952036352
$ws.Range("Q4").Value = 0.20252694275
$ws.Range("R4").Value = 1.82274248475
$ws.Range("S4").Value = 0.000848525125584423
$ws.Range("T4").Value = 0.0008485251255844231
$ws.Range("I5").Value = 0.07940559828131502
$ws.Range("J5").Value = 0.07940559828131502
$ws.Range("M5").Value = 1.937269333333333
$ws.Range("N5").Value = 5.811808
$ws.Range("O5").Value = 0.5832046952539398
$ws.Range("P5").Value = 0.5832046952539398
$ws.Range("Q5").Value = 11.05325614073956
$ws.Range("R5").Value = 99.479305266656
$ws.Range("S5").Value = 0.0463097177471111
$ws.Range("T5").Value = 0.0463097177471111
$ws.Range("I6").Value = 0.07940559828131502
$ws.Range("J6").Value = 0.07940559828131502
$ws.Range("M6").Value = 0.9034129999999999
$ws.Range("N6").Value = 2.710239
$ws.Range("O6").Value = 0.2719677095424251
$ws.Range("P6").Value = 0.2719677095424251
$ws.Range("Q6").Value = 5.154500263880333
$ws.Range("R6").Value = 46.390502374923
$ws.Range("S6").Value = 0.02159575868941517
$ws.Range("T6").Value = 0.02159575868941517
$ws.Range("I7").Value = 0.07940559828131502
$ws.Range("J7").Value = 0.07940559828131502
$ws.Range("M7").Value = 0.4810833333333333
$ws.Range("N7").Value = 1.44325
$ws.Range("O7").Value = 0.1448275952036352
$ws.Range("P7").Value = 0.1448275952036352
$ws.Range("Q7").Value = 2.744862171138889
$ws.Range("R7").Value = 24.70375954025
$ws.Range("S7").Value = 0.01150012184478876
$ws.Range("T7").Value = 0.01150012184478876
$ws.Range("G8").Value = 37.57387866666667
$ws.Range("H8").Value = 112.721636
$ws.Range("I8").Value = 0.5229220082886389
$ws.Range("J8").Value = 0.5229220082886389
$ws.Range("M8").Value = 1.937269333333333
$ws.Range("N8").Value = 5.811808
$ws.Range("O8").Value = 0.5832046952539398
$ws.Range("P8").Value = 0.5832046952539398
$ws.Range("Q8").Value = 72.79072287532088
$ws.Range("R8").Value = 655.1165058778879
$ws.Range("S8").Value = 0.3049705704855538
$ws.Range("T8").Value = 0.3049705704855538
$ws.Range("G9").Value = 37.57387866666667
$ws.Range("H9").Value = 112.721636
$ws.Range("I9").Value = 0.5229220082886389
$ws.Range("J9").Value = 0.5229220082886389
$ws.Range("M9").Value = 0.9034129999999999
$ws.Range("N9").Value = 2.710239
$ws.Range("O9").Value = 0.2719677095424251
$ws.Range("P9").Value = 0.2719677095424251
$ws.Range("Q9").Value = 33.94473044788933
$ws.Range("R9").Value = 305.5025740310039
$ws.Range("S9").Value = 0.1422179008635861
$ws.Range("T9").Value = 0.1422179008635861
$ws.Range("G10").Value = 37.57387866666667
$ws.Range("H10").Value = 112.721636
$ws.Range("I10").Value = 0.5229220082886389
$ws.Range("J10").Value = 0.5229220082886389
$ws.Range("M10").Value = 0.4810833333333333
$ws.Range("N10").Value = 1.44325
$ws.Range("O10").Value = 0.1448275952036352
$ws.Range("P10").Value = 0.1448275952036352
$ws.Range("Q10").Value = 18.07616679522222
$ws.Range("R10").Value = 162.685501157
$ws.Range("S10").Value = 0.07573353693949894
$ws.Range("T10").Value = 0.07573353693949894
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.764088
$ws.Range("H11").Value = 5.292263999999999
$ws.Range("I11").Value = 0.02455111030568848
$ws.Range("J11").Value = 0.02455111030568848
$ws.Range("M11").Value = 1.937269333333333
$ws.Range("N11").Value = 5.811808
$ws.Range("O11").Value = 0.5832046952539398
$ws.Range("P11").Value = 0.5832046952539398
$ws.Range("Q11").Value = 3.417513583701333
$ws.Range("R11").Value = 30.757622253312
$ws.Range("S11").Value = 0.01431832280397491
$ws.Range("T11").Value = 0.01431832280397491
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.764088
$ws.Range("H12").Value = 5.292263999999999
$ws.Range("I12").Value = 0.02455111030568848
$ws.Range("J12").Value = 0.02455111030568848
$ws.Range("M12").Value = 0.9034129999999999
$ws.Range("N12").Value = 2.710239
$ws.Range("O12").Value = 0.2719677095424251
$ws.Range("P12").Value = 0.2719677095424251
$ws.Range("Q12").Value = 1.593700032344
$ws.Range("R12").Value = 14.343300291096
$ws.Range("S12").Value = 0.006677109236561522
$ws.Range("T12").Value = 0.006677109236561522
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.764088
$ws.Range("H13").Value = 5.292263999999999
$ws.Range("I13").Value = 0.02455111030568848
$ws.Range("J13").Value = 0.02455111030568848
$ws.Range("M13").Value = 0.4810833333333333
$ws.Range("N13").Value = 1.44325
$ws.Range("O13").Value = 0.1448275952036352
$ws.Range("P13").Value = 0.1448275952036352
$ws.Range("Q13").Value = 0.8486733353333332
$ws.Range("R13").Value = 7.638060017999999
$ws.Range("S13").Value = 0.003555678265152046
$ws.Range("T13").Value = 0.003555678265152046
$ws.Range("G14").Value = 26.389162
$ws.Range("H14").Value = 79.167486
$ws.Range("I14").Value = 0.3672624195259435
$ws.Range("J14").Value = 0.3672624195259436
$ws.Range("M14").Value = 1.937269333333333
$ws.Range("N14").Value = 5.811808
$ws.Range("O14").Value = 0.5832046952539398
$ws.Range("P14").Value = 0.5832046952539398
$ws.Range("Q14").Value = 51.12291427496533
$ws.Range("R14").Value = 460.106228474688
$ws.Range("S14").Value = 0.2141891674578525
$ws.Range("T14").Value = 0.2141891674578525
$ws.Range("G15").Value = 26.389162
$ws.Range("H15").Value = 79.167486
$ws.Range("I15").Value = 0.3672624195259435
$ws.Range("J15").Value = 0.3672624195259436
$ws.Range("M15").Value = 0.9034129999999999
$ws.Range("N15").Value = 2.710239
$ws.Range("O15").Value = 0.2719677095424251
$ws.Range("P15").Value = 0.2719677095424251
$ws.Range("Q15").Value = 23.840312009906
$ws.Range("R15").Value = 214.5628080891539
$ws.Range("S15").Value = 0.09988351903948006
$ws.Range("T15").Value = 0.09988351903948008
$ws.Range("G16").Value = 26.389162
$ws.Range("H16").Value = 79.167486
$ws.Range("I16").Value = 0.3672624195259435
$ws.Range("J16").Value = 0.3672624195259436
$ws.Range("M16").Value = 0.4810833333333333
$ws.Range("N16").Value = 1.44325
$ws.Range("O16").Value = 0.1448275952036352
$ws.Range("P16").Value = 0.1448275952036352
$ws.Range("Q16").Value = 12.69538601883333
$ws.Range("R16").Value = 114.2584741695
$ws.Range("S16").Value = 0.05318973302861098
$ws.Range("T16").Value = 0.05318973302861099
